# Update countries & provincias Spain
#
# The source data feed (COVID-19 daily snapshot) refreshed the figures for a
# handful of countries. The sheet is kept sorted by "Casos totales" (column B)
# descending, so after pushing in the new numbers we re-sort the whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-country totals: Pais, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @(
    @("Israel",      13107, 125, 3247, 9702, 167, 7, 158),
    @("Polonia",       8379,  0,  981, 7066, 160, 0, 332),
    @("Noruega",       6992, 55,   32, 6798,  63, 1, 162),
    @("Chequia",       6553,  4, 1183, 5194,  86, 3, 176),
    @("Armenia",       1248, 47,  523,  705,  30, 1,  20),
    @("Lituania",      1239, 90,  228,  978,  14, 0,  33),
    @("Letonia",        712, 30,   88,  619,   5, 0,   5),
    @("Taiwan",         398,  3,  178,  214,   0, 0,   6),
    @("Montenegro",     305,  2,   55,  245,   7, 0,   5),
    @("Guatemala",      235, 21,   21,  207,   3, 0,   7),
    @("Islas Feroe",    184,  0,  173,   11,   0, 0,   0)
)

$countryRange = $ws.Range("A4:A216")

foreach ($row in $updates) {
    $name = $row[0]
    $found = $countryRange.Find($name, [Type]::Missing, [Type]::Missing, 1)
    $r = $found.Row

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Re-sort the data block by Casos totales (column B) descending, keeping it
# in the same order the site publishes the ranking.
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2, [Type]::Missing, [Type]::Missing, 1, [Type]::Missing, [Type]::Missing, 1)

# Letonia's updated total (712) now ties with Crucero's; the published
# ranking lists Letonia ahead of Crucero for that tie, so swap the pair back
# if the stable sort put them in the other order.
$letoniaCell = $ws.Range("A4:A216").Find("Letonia", [Type]::Missing, [Type]::Missing, 1)
$cruceroCell = $ws.Range("A4:A216").Find("Crucero", [Type]::Missing, [Type]::Missing, 1)
if ($letoniaCell.Row -gt $cruceroCell.Row) {
    $r1 = $cruceroCell.Row
    $r2 = $letoniaCell.Row
    for ($c = 1; $c -le 8; $c++) {
        $tmp = $ws.Cells.Item($r1, $c).Value2
        $ws.Cells.Item($r1, $c).Value = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r2, $c).Value = $tmp
    }
}
